$wb = $excel.ActiveWorkbook

# --- Sheet "Clients": just a view/selection change (Orders stays the active tab) ---
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Range("G4").Select() | Out-Null

# --- Sheet "Orders": data + formula fixes ---
$wsOrders = $wb.Worksheets.Item("Orders")
$wsOrders.Activate() | Out-Null

# Fix the AddressId column (C2:C8): was a literal "NULL" text, now real numeric ids
$wsOrders.Range("C2").Value = 2
$wsOrders.Range("C3").Value = 3
$wsOrders.Range("C4").Value = 4
$wsOrders.Range("C5").Value = 2
$wsOrders.Range("C6").Value = 3
$wsOrders.Range("C7").Value = 4
$wsOrders.Range("C8").Value = 1

# Rename the "Date" header (E26) of the ProductionOrder insert table to "ProductionDate"
$wsOrders.Range("E26").Value = "ProductionDate"

# Fix the ProductionOrder INSERT formula so the ProductId is quoted as a string literal
$wsOrders.Range("G27").Formula = '="INSERT INTO ProductionOrder(" & $A$26 & "," & $B$26 & "," & $C$26 & "," & $D$26 & "," & $E$26 & ") VALUES(" & A27 & "," & B27 & ",''" & C27 & "'', " & D27 & ", TO_DATE(" & TEXT(E27,"''dd/MM/AAAA''") & ", ''dd/MM/YYYY''));"'
$wsOrders.Range("G28:G39").Formula = '="INSERT INTO ProductionOrder(" & $A$26 & "," & $B$26 & "," & $C$26 & "," & $D$26 & "," & $E$26 & ") VALUES(" & A28 & "," & B28 & ",''" & C28 & "'', " & D28 & ", TO_DATE(" & TEXT(E28,"''dd/MM/AAAA''") & ", ''dd/MM/YYYY''));"'

# Column C got a bit wider to fit the new numeric ids
$wsOrders.Columns("C").ColumnWidth = 10.75
